# feat: Ajout des préconisations d'heures dans l'import Excel
#
# Adds three new header columns ("preco CM", "Preco TD", "Préco TP") to the
# "Ressources" sheet, then leaves the "SAE" sheet as the active/selected tab
# (mirroring the view-state changes captured in the workbook diff).

$wb  = $excel.ActiveWorkbook
$wsRessources = $wb.Worksheets.Item("Ressources")
$wsSAE        = $wb.Worksheets.Item("SAE")

# New header cells on the "Ressources" sheet (columns S, T, U of row 1).
$wsRessources.Range("S1").Value = "preco CM"
$wsRessources.Range("T1").Value = "Preco TD"
$wsRessources.Range("U1").Value = "Préco TP"

# Reproduce the resulting view state on "Ressources": scrolled right so
# column H is the left-most visible column, with cell V1 selected.
$wsRessources.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$wsRessources.Range("V1").Select()

# The "SAE" sheet ends up as the active/selected tab in the saved workbook.
$wsSAE.Activate()
